$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.056748553031213
$ws.Range("D2").Value = 1.054271244083026
$ws.Range("E2").Value = 1.062098227298203
$ws.Range("F2").Value = 1.07078119812752
$ws.Range("I2").Value = 1.041198566114515
$ws.Range("J2").Value = 1.061748569773783
$ws.Range("K2").Value = 1.057014914388477
$ws.Range("L2").Value = 1.064820496232768
$ws.Range("M2").Value = 1.073480121093358
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.058371590284791
$ws.Range("D3").Value = 1.055503062703236
$ws.Range("E3").Value = 1.063544809407952
$ws.Range("F3").Value = 1.072313733926124
$ws.Range("I3").Value = 1.041579471626974
$ws.Range("J3").Value = 1.063020560694667
$ws.Range("K3").Value = 1.058058873017663
$ws.Range("L3").Value = 1.066080238859819
$ws.Range("M3").Value = 1.074827314628468
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059419962158909
$ws.Range("D4").Value = 1.056298266789985
$ws.Range("E4").Value = 1.064479359051929
$ws.Range("F4").Value = 1.073303973633466
$ws.Range("I4").Value = 1.041823749525428
$ws.Range("J4").Value = 1.063841422877486
$ws.Range("K4").Value = 1.058731960273043
$ws.Range("L4").Value = 1.066893373492231
$ws.Range("M4").Value = 1.07569712201977
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.059860266817199
$ws.Range("D5").Value = 1.056632130791777
$ws.Range("E5").Value = 1.064871896053509
$ws.Range("F5").Value = 1.073719940487427
$ws.Range("I5").Value = 1.041925921480718
$ws.Range("J5").Value = 1.064185994068419
$ws.Range("K5").Value = 1.059014352150663
$ws.Range("L5").Value = 1.067234743140964
$ws.Range("M5").Value = 1.076062337903094
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.059934170930517
$ws.Range("D6").Value = 1.05668816242983
$ws.Range("E6").Value = 1.064937784587705
$ws.Range("F6").Value = 1.073789764050707
$ws.Range("I6").Value = 1.04194304603983
$ws.Range("J6").Value = 1.064243818858654
$ws.Range("K6").Value = 1.05906173351047
$ws.Range("L6").Value = 1.067292033133768
$ws.Range("M6").Value = 1.076123633026364
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059425847210402
$ws.Range("D7").Value = 1.056302729615854
$ws.Range("E7").Value = 1.064484605505458
$ws.Range("F7").Value = 1.073309533087467
$ws.Range("I7").Value = 1.041825116801328
$ws.Range("J7").Value = 1.063846029084298
$ws.Range("K7").Value = 1.058735735853815
$ws.Range("L7").Value = 1.066897936735087
$ws.Range("M7").Value = 1.075702003814111
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.05729745348192
$ws.Range("D8").Value = 1.054687933129256
$ws.Range("E8").Value = 1.062587418960485
$ws.Range("F8").Value = 1.071299423544587
$ws.Range("I8").Value = 1.041327750677337
$ws.Range("J8").Value = 1.062178905413879
$ws.Range("K8").Value = 1.057368230733854
$ws.Range("L8").Value = 1.065246651367257
$ws.Range("M8").Value = 1.073935813056145
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.053532394482878
$ws.Range("D9").Value = 1.051827871752655
$ws.Range("E9").Value = 1.059232594492013
$ws.Range("F9").Value = 1.067746147879224
$ws.Range("I9").Value = 1.040434409813154
$ws.Range("J9").Value = 1.05922401904455
$ws.Range("K9").Value = 1.054939666980677
$ws.Range("L9").Value = 1.062321201661764
$ws.Range("M9").Value = 1.070808528432474
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051011933803925
$ws.Range("D10").Value = 1.049910935958148
$ws.Range("E10").Value = 1.056987645181501
$ws.Range("F10").Value = 1.065369228472952
$ws.Range("I10").Value = 1.039827306375655
$ws.Range("J10").Value = 1.057242045399173
$ws.Range("K10").Value = 1.053307557728329
$ws.Range("L10").Value = 1.060359895997202
$ws.Range("M10").Value = 1.06871307475096
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.0499179329975
$ws.Range("D11").Value = 1.049078360654492
$ws.Range("E11").Value = 1.056013451549505
$ws.Range("F11").Value = 1.064337962524581
$ws.Range("I11").Value = 1.039561648224685
$ws.Range("J11").Value = 1.056380862518024
$ws.Range("K11").Value = 1.05259764818152
$ws.Range("L11").Value = 1.059507913780696
$ws.Range("M11").Value = 1.06780309656176
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049511165769607
$ws.Range("D12").Value = 1.048768716818578
$ws.Range("E12").Value = 1.05565126492538
$ws.Range("F12").Value = 1.063954586920673
$ws.Range("I12").Value = 1.039462550214834
$ws.Range("J12").Value = 1.056060524425825
$ws.Range("K12").Value = 1.05233346832482
$ws.Range("L12").Value = 1.059191031299627
$ws.Range("M12").Value = 1.067464685007011
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049598437319535
$ws.Range("D13").Value = 1.048835154177492
$ws.Range("E13").Value = 1.055728970175659
$ws.Range("F13").Value = 1.064036836825271
$ws.Range("I13").Value = 1.039483826167673
$ws.Range("J13").Value = 1.056129258840631
$ws.Range("K13").Value = 1.052390158002105
$ws.Range("L13").Value = 1.059259022745975
$ws.Range("M13").Value = 1.067537293922256
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049884317857896
$ws.Range("D14").Value = 1.049052773367732
$ws.Range("E14").Value = 1.055983519833205
$ws.Range("F14").Value = 1.064306279115253
$ws.Range("I14").Value = 1.039553465354696
$ws.Range("J14").Value = 1.056354392637225
$ws.Range("K14").Value = 1.052575820989136
$ws.Range("L14").Value = 1.059481728763213
$ws.Range("M14").Value = 1.067775131677464
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050060404042584
$ws.Range("D15").Value = 1.049186804031039
$ws.Range("E15").Value = 1.056140312553987
$ws.Range("F15").Value = 1.064472249016427
$ws.Range("I15").Value = 1.039596316503418
$ws.Range("J15").Value = 1.056493044235097
$ws.Range("K15").Value = 1.052690149205888
$ws.Range("L15").Value = 1.059618889581003
$ws.Range("M15").Value = 1.067921617395446
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.051084482688779
$ws.Range("D16").Value = 1.049966137241743
$ws.Range("E16").Value = 1.057052253787819
$ws.Range("F16").Value = 1.065437626243748
$ws.Range("I16").Value = 1.039844878412156
$ws.Range("J16").Value = 1.057299135783164
$ws.Range("K16").Value = 1.053354604101427
$ws.Range("L16").Value = 1.060416381085135
$ws.Range("M16").Value = 1.068773410718816
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05172614882067
$ws.Range("D17").Value = 1.05045430935014
$ws.Range("E17").Value = 1.057623717044401
$ws.Range("F17").Value = 1.066042627279511
$ws.Range("I17").Value = 1.040000048546624
$ws.Range("J17").Value = 1.057803972532112
$ws.Range("K17").Value = 1.053770538129589
$ws.Range("L17").Value = 1.060915890757344
$ws.Range("M17").Value = 1.069307006528784
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.052100169409515
$ws.Range("D18").Value = 1.050738808389613
$ws.Range("E18").Value = 1.057956838424081
$ws.Range("F18").Value = 1.066395317623893
$ws.Range("I18").Value = 1.040090288816416
$ws.Range("J18").Value = 1.058098148843666
$ws.Range("K18").Value = 1.054012837622916
$ws.Range("L18").Value = 1.061206984099435
$ws.Range("M18").Value = 1.069617990550154
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.052227658321592
$ws.Range("D19").Value = 1.050835774158416
$ws.Range("E19").Value = 1.058070389990804
$ws.Range("F19").Value = 1.066515542896171
$ws.Range("I19").Value = 1.040121013085403
$ws.Range("J19").Value = 1.058198407065318
$ws.Range("K19").Value = 1.054095403515765
$ws.Range("L19").Value = 1.061306195207369
$ws.Range("M19").Value = 1.069723985374168
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051657330314201
$ws.Range("D20").Value = 1.050401958348365
$ws.Range("E20").Value = 1.057562425562437
$ws.Range("F20").Value = 1.065977736796727
$ws.Range("I20").Value = 1.039983427988203
$ws.Range("J20").Value = 1.057749838021015
$ws.Range("K20").Value = 1.053725944236867
$ws.Range("L20").Value = 1.060862325283221
$ws.Range("M20").Value = 1.069249783026547
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049800144488591
$ws.Range("D21").Value = 1.048988700725993
$ws.Range("E21").Value = 1.055908570430239
$ws.Range("F21").Value = 1.064226943925456
$ws.Range("I21").Value = 1.039532969991311
$ws.Range("J21").Value = 1.056288109021221
$ws.Range("K21").Value = 1.052521161396494
$ws.Range("L21").Value = 1.059416159014081
$ws.Range("M21").Value = 1.067705105647619
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048630099317789
$ws.Range("D22").Value = 1.048097878839509
$ws.Range("E22").Value = 1.054866824953453
$ws.Range("F22").Value = 1.063124308702749
$ws.Range("I22").Value = 1.039247313202295
$ws.Range("J22").Value = 1.055366416507998
$ws.Range("K22").Value = 1.051760840205648
$ws.Range("L22").Value = 1.058504471856484
$ws.Range("M22").Value = 1.066731556860409
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049250590110704
$ws.Range("D23").Value = 1.048570336502217
$ws.Range("E23").Value = 1.055419257460372
$ws.Range("F23").Value = 1.063709014617792
$ws.Range("I23").Value = 1.039398977291677
$ws.Range("J23").Value = 1.055855277080088
$ws.Range("K23").Value = 1.052164171373781
$ws.Range("L23").Value = 1.058988007396317
$ws.Range("M23").Value = 1.067247879350655
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.051688427215024
$ws.Range("D24").Value = 1.050425614265179
$ws.Range("E24").Value = 1.057590121173728
$ws.Range("F24").Value = 1.066007058620116
$ws.Range("I24").Value = 1.039990938930814
$ws.Range("J24").Value = 1.057774299961525
$ws.Range("K24").Value = 1.053746095248628
$ws.Range("L24").Value = 1.060886530025574
$ws.Range("M24").Value = 1.06927564064713
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.054507542535221
$ws.Range("D25").Value = 1.052569037721631
$ws.Range("E25").Value = 1.060101340116731
$ws.Range("F25").Value = 1.068666137540185
$ws.Range("I25").Value = 1.040667381931662
$ws.Range("J25").Value = 1.059990016044024
$ws.Range("K25").Value = 1.055569782746296
$ws.Range("L25").Value = 1.063079405866312
$ws.Range("M25").Value = 1.071618838030408
